$wb = $excel.ActiveWorkbook

# --- English sheet: update selection (I10 -> B6) ---
$wsEnglish = $wb.Worksheets.Item("English")
$wsEnglish.Range("B6").Select()

# --- Kannada sheet: rewrite the long "how to protect yourself" answer as an <ul><li> list ---
$wsKannada = $wb.Worksheets.Item("Kannada")
$wsKannada.Range("B6").Value = '<ul><li>ನಿಯಮಿತವಾಗಿ ಮತ್ತು ಸಂಪೂರ್ಣವಾಗಿ ನಿಮ್ಮ ಕೈಗಳನ್ನು ಆಲ್ಕೋಹಾಲ್ ಆಧಾರಿತ ಹ್ಯಾಂಡ್ ರಬ್ನಿಂದ ಸ್ವಚ್ clean ಗೊಳಿಸಿ ಅಥವಾ ಸೋಪ್ ಮತ್ತು ನೀರಿನಿಂದ ತೊಳೆಯಿರಿ.</li>
<li>ಏಕೆ? ನಿಮ್ಮ ಕೈಗಳನ್ನು ಸೋಪ್ ಮತ್ತು ನೀರಿನಿಂದ ತೊಳೆಯುವುದು ಅಥವಾ ಆಲ್ಕೋಹಾಲ್ ಆಧಾರಿತ ಹ್ಯಾಂಡ್ ರಬ್ ಬಳಸುವುದರಿಂದ ನಿಮ್ಮ ಕೈಯಲ್ಲಿರುವ ವೈರಸ್‌ಗಳನ್ನು ಕೊಲ್ಲುತ್ತದೆ.</li>
<li>ನಿಮ್ಮ ಮತ್ತು ಕೆಮ್ಮು ಅಥವಾ ಸೀನುವವರ ನಡುವೆ ಕನಿಷ್ಠ 1 ಮೀಟರ್ (3 ಅಡಿ) ಅಂತರವನ್ನು ಕಾಪಾಡಿಕೊಳ್ಳಿ.</li>
<li>ಏಕೆ? ಯಾರಾದರೂ ಕೆಮ್ಮಿದಾಗ ಅಥವಾ ಸೀನುವಾಗ ಅವರು ಮೂಗು ಅಥವಾ ಬಾಯಿಯಿಂದ ಸಣ್ಣ ದ್ರವ ಹನಿಗಳನ್ನು ಸಿಂಪಡಿಸುತ್ತಾರೆ, ಅದು ವೈರಸ್ ಅನ್ನು ಹೊಂದಿರಬಹುದು. ನೀವು ತುಂಬಾ ಹತ್ತಿರದಲ್ಲಿದ್ದರೆ, ಕೆಮ್ಮುವ ವ್ಯಕ್ತಿಗೆ ಕಾಯಿಲೆ ಇದ್ದರೆ ನೀವು Covid 19 ವೈರಸ್ ಸೇರಿದಂತೆ ಹನಿಗಳಲ್ಲಿ ಉಸಿರಾಡಬಹುದು.</li>
<li>ಕಣ್ಣು, ಮೂಗು ಮತ್ತು ಬಾಯಿಯನ್ನು ಸ್ಪರ್ಶಿಸುವುದನ್ನು ತಪ್ಪಿಸಿ.</li><li>ನೀವು ಕೆಮ್ಮುವಾಗ ಅಥವಾ ಸೀನುವಾಗ ನಿಮ್ಮ ಬಗ್ಗಿದ  ಮೊಣಕೈ ಅಥವಾ ಟಿಶ್ಯೂ ದಿಂದ  ನಿಮ್ಮ ಬಾಯಿ ಮತ್ತು ಮೂಗನ್ನು ಮುಚ್ಚಿ. ನಂತರ ಬಳಸಿ ಟಿಸ್ಸುದನ್ನು  ತಕ್ಷಣ ವಿಲೇವಾರಿ ಮಾಡಿ.</li><li>ನಿಮಗೆ ಅನಾರೋಗ್ಯ ಅನಿಸಿದರೆ ಮನೆಯಲ್ಲೇ ಇರಿ. ನಿಮಗೆ ಜ್ವರ, ಕೆಮ್ಮು ಮತ್ತು ಉಸಿರಾಟದ ತೊಂದರೆ ಇದ್ದರೆ, ವೈದ್ಯಕೀಯ ಚಿಕಿತ್ಸೆ ಪಡೆಯಿರಿ ಮತ್ತು ಮುಂಚಿತವಾಗಿ ಕರೆ ಮಾಡಿ.</li><li>ನಿಮ್ಮ ಸ್ಥಳೀಯ ಆರೋಗ್ಯ ಪ್ರಾಧಿಕಾರದ ನಿರ್ದೇಶನಗಳನ್ನು ಅನುಸರಿಸಿ</li></ul>'
$wsKannada.Rows.Item(6).RowHeight = 204

# --- Kannada sheet: add new Q&A row 7 ---
$wsKannada.Range("A7").Value = 'ನಾನು ಕರೋನವೈರಸ್ ಲಕ್ಷಣಗಳನ್ನು ಹೊಂದಿದ್ದರೆ ಏನು ಮಾಡಬೇಕು?'
$wsKannada.Range("B7").Value = '1075 (ಟೋಲ್ ಫ್ರೀ) ಅಥವಾ + 91-11-23978046 ಗೆ ಫೋನ್ ಕರೆ ಮಾಡಿ'
$wsKannada.Range("B8").Select()

# --- Add new "Bengali" sheet after Kannada (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBengali = $wb.Worksheets.Add($null, $lastSheet)
$wsBengali.Name = "Bengali"
$wsBengali.Columns.Item(1).ColumnWidth = 65.5

$wsBengali.Range("A1").Value = "question"
$wsBengali.Range("B1").Value = "answer"

# Column A questions first (rows 3-6), matching shared-string creation order
$wsBengali.Range("A3").Value = 'করোনা ভাইরাসের লক্ষণগুলি কি কি?'
$wsBengali.Range("A4").Value = 'করোনা ভাইরাস কীভাবে ছড়িয়ে পড়ে?'
$wsBengali.Range("A5").Value = 'আমি নিজেকে কীভাবে করোনা ভাইরাস থেকে রক্ষা করব?'
$wsBengali.Range("A6").Value = 'আপনার  করোনা ভাইরাসের  লক্ষণ থাকলে কী করবেন?'

# Then column B answers: B6, B3, B4, B5
$wsBengali.Range("B6").Value = '১০৭৫ (টোল ফ্রি) ও +৯১-১১-২৩৯৭৮০৪৬ নম্বর এ কল করুন'
$wsBengali.Range("B3").Value = 'কোভিড ১৯ এর সবচেয়ে সাধারণ লক্ষণগুলি হল জ্বর, ক্লান্তি এবং শুকনো কাশি। কিছু রোগীর শরীরে ব্যথা, অবরুদ্ধ নাক, সর্দি, গলা ব্যথা বা ডায়রিয়া হতে পারে। লক্ষণগুলি সাধারণত হালকা এবং ধীরে ধীরে শুরু হয়। কিছু লোকের ক্ষেত্রে সংক্রামিত হলেও কোনও উপসর্গ বিকাশ করে না এবং অসুস্থ বোধ করে না। প্রায় ৫ জনের মধ্যে ১ জন বিশেষ চিকিত্সার প্রয়োজন ছাড়াই এই রোগ থেকে সেরে উঠেন। কোভিড আক্রান্তের মধ্যে প্রতি ৬ জনের মধ্যে ১ জন গুরুতর অসুস্থ হয়ে পড়ে এবং শ্বাস নিতে অসুবিধা বোধ করে। বয়স্ক ব্যক্তিরা এবং যারা উচ্চ রক্তচাপ, হার্টের সমস্যা বা ডায়াবেটিসের মতো রোগে চিকিত্সাধীন তাদের গুরুতর অসুস্থ হওয়ার সম্ভাবনা বেশি থাকে। জ্বর, কাশি এবং শ্বাস নিতে অসুবিধাজনিত ব্যক্তিদের চিকিত্সার পরামর্শ নেওয়া উচিত।'
$wsBengali.Range("B4").Value = 'লোকেরা ভাইরাসযুক্ত অন্যদের কাছ থেকে কোভিড ১৯ সংক্রামিত হতে পারে। কোভিড ১৯ সংক্রামিত কোনো ব্যক্তির কাশি বা শ্বাস ছাড়ার সময় নাক বা মুখ থেকে ছোট ফোঁটাগুলির মাধ্যমে এই রোগটি এক ব্যক্তি থেকে অন্য ব্যক্তির কাছে ছড়িয়ে পড়ে। এই ফোঁটাগুলি ব্যক্তিটির চারপাশে বস্তু এবং পৃষ্ঠগুলিতে অবতরণ করে। অন্যান্য ব্যক্তিরা তখন এই জিনিসগুলি বা পৃষ্ঠের স্পর্শ করে, তারপর তাদের চোখ, নাক বা মুখ স্পর্শ করে তখন তারা সংক্রমিত হন। যদি কোনো বেক্তি, কোভিড ১৯ সংক্রামিত কোনো ব্যক্তির কাশি বা শ্বাস ছাড়ার সময়ের ফোঁটাগুলি শ্বাসের মাধ্যমে শরীরের ভিতরে প্রবেশ করান তাহলেও সংক্রমন হতে পারে। এজন্য অসুস্থ ব্যক্তির থেকে ১ মিটার বা ৩ ফুট বেশি দূরে থাকা গুরুত্বপূর্ণ।'
$wsBengali.Range("B5").Value = '<ul><li>নিয়মিত এবং সঠিকভাবে সাবান এবং জল দিয়ে আপনার হাত ধুলে, বা অ্যালকোহল ভিত্তিক হাতের স্যানিটাইজার দিয়ে পরিষ্কার করলে আপনার হাতের ভাইরাসগুলি মরে যাবে। </li> <li>নিজের এবং কাশি বা হাঁচি হয় এমন কারও মধ্যে অন্তত ১ মিটার বা ৩ ফুট দূরত্ব বজায় রাখুন। যখন কেউ কাশি বা হাঁচে তখন তারা তাদের নাক বা মুখ থেকে ছোট তরল ফোঁটা স্প্রে করে যার মধ্যে ভাইরাস থাকতে পারে। যদি আপনি কোনো কোভিড ১৯ সংক্রমিত ব্যক্তির খুব কাছাকাছি থাকেন তাহলে তার কাশির সাথে ভাইরাস আপনার শরীরের ভিতরে প্রবেশ করতে পারে। </li> <li>চোখ, নাক এবং মুখ স্পর্শ করা এড়িয়ে চলুন।</li> <li>কাশি বা হাঁচি হলে আপনার বাঁকানো কনুই বা টিস্যু দিয়ে আপনার মুখ এবং নাকটি ঢেকে রাখুন। তারপরে অবিলম্বে ব্যবহৃত টিস্যুগুলি ফেলে দিন। </li> <li>আপনি অসুস্থ বোধ করলে বাড়িতেই থাকুন। আপনার যদি জ্বর, কাশি এবং শ্বাস নিতে সমস্যা হয় তবে চিকিত্সার পরামর্শ নিন এবং আগে থেকেই কল করুন।</li> <li>আপনার স্থানীয় স্বাস্থ্য কর্তৃপক্ষের নির্দেশাবলী অনুসরণ করুন।</li></ul>'

# Finally the intro "what is coronavirus" question/answer in row 2
$wsBengali.Range("A2").Value = 'করোনা ভাইরাস বা কোভিড ১৯ কি?'
$wsBengali.Range("B2").Value = 'করোনা ভাইরাস হল এক ধরণের ভাইরাস যা প্রাণী বা মানুষের সাধারণ সর্দি বা আরও গুরুতর অসুস্থতার কারণ হতে পারে।কোভিড ১৯ হল এক ধরণের করোনা  ভাইরাস যা সম্প্রতি ২০১৯ সালে চিনে পাওয়া গিয়েছিল এবং এটি ভারত সহ বিশ্বের বিভিন্ন জায়গায় ছড়িয়ে পড়েছে।'

$wsBengali.Range("J9").Select()
$wsBengali.Activate()
